$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.759.30"
$ws.Range("E2").Value = "  +3.53%  "
$ws.Range("D3").Value = "3.464.83"
$ws.Range("E3").Value = "  +2.24%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.63"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.68"
$ws.Range("E6").Value = "  +7.82%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "3.466.51"
$ws.Range("E8").Value = "  +2.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.563"
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.25"
$ws.Range("E10").Value = "  +1.44%  "
$ws.Range("E11").Value = "  +4.12%  "
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").Value = "4.059.05"
$ws.Range("E13").Value = "  +1.83%  "
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.62"
$ws.Range("E15").Value = "  +2.65%  "
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("D17").Value = "65.697.94"
$ws.Range("E17").Value = "  +3.20%  "
$ws.Range("D18").Value = "3.455.75"
$ws.Range("E18").Value = "  +1.59%  "
$ws.Range("E19").Value = "  +2.37%  "
$ws.Range("E20").Value = "  +2.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.95"
$ws.Range("E21").Value = "  +2.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.93"
$ws.Range("E22").Value = "  +2.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.85"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("E26").Value = "  +3.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.96"
$ws.Range("E27").Value = "  +3.63%  "
$ws.Range("E28").Value = "  +2.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.982"
$ws.Range("E29").Value = "  -1.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.33"
$ws.Range("E30").Value = "  +6.68%  "
$ws.Range("E31").Value = "  +4.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.07"
$ws.Range("E32").Value = "  +4.99%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.25"
$ws.Range("E33").Value = "  +2.26%  "
$ws.Range("E34").Value = "  +6.12%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.83"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.899"
$ws.Range("E38").Value = "  +10.58%  "
$ws.Range("E39").Value = "  +3.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0742"
$ws.Range("E40").Value = "  +2.37%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.821.26"
$ws.Range("E41").Value = "  +1.85%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.64"
$ws.Range("E42").Value = "  +5.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.16"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "43.14"
$ws.Range("E44").Value = "  +1.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "26.51"
$ws.Range("E45").Value = "  +4.74%  "
$ws.Range("E46").Value = "  +1.48%  "
$ws.Range("E47").Value = "  +2.76%  "
$ws.Range("E48").Value = "  +7.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "346.55"
$ws.Range("E49").Value = "  +6.80%  "
$ws.Range("E50").Value = "  +4.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "32.62"
$ws.Range("E51").Value = "  +9.37%  "
